$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "scenarios"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "present"

$ws2.Range("A1").Value = "scenarios"
$ws2.Range("B1").Value = "years"
$ws2.Range("C1").Value = "boiler oil"
$ws2.Range("D1").Value = "boiler gas"
$ws2.Range("E1").Value = "air-water heat pump"
$ws2.Range("A2").Value = "BAU"
$ws2.Range("B2").Value = 2015
$ws2.Range("C2").Value = 1000
$ws2.Range("D2").Value = 2000
$ws2.Range("E2").Value = 100

$ws1.Range("A1:F2").Select() | Out-Null
$ws2.Range("F4").Select() | Out-Null
